$d = $word.ActiveDocument

# 1. Add " – 25/09/2019" (bold) after "9" in the title "Minutes (Week 9)"
$r = $d.Content
$r.Find.Execute("Minutes (Week 9", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Collapse(0)
    $r.Font.Bold = $true
    $r.InsertAfter(" - 25/09/2019")
}

Write-Output "done"
